$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = 5.898099999999991
$ws.Range("D4").Value = -7.134799999999996
$ws.Range("D6").Value = -7.983899999999998
$ws.Range("A9").Value = -20.35659999999998
$ws.Range("D10").Value = -7.689699999999993
$ws.Range("B11").Value = 5.813399999999999
$ws.Range("D11").Value = -8.123699999999999
$ws.Range("E12").Value = 12.16899999999999
$ws.Range("E17").Value = 13.3503
$ws.Range("A18").Value = -22.87180000000002
$ws.Range("E19").Value = 12.95459999999999
$ws.Range("A20").Value = -22.13470000000002
$ws.Range("C21").Value = -12.6984
$ws.Range("D21").Value = -7.846699999999999
